# Actualización automática de tasas-transfi.xlsx
# Updates the "Conversión del día" note on Hoja1 and the N10/O10/N12/O12
# rate cells on the "tasas" sheet.

$wb = $excel.ActiveWorkbook

# --- Hoja1: update the conversion text cell (A1) ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")

$newText = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 2.35 = 8894.52 pesos`n✅ 8894.52 pesos = 2.34 = 936.5 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

$wsHoja1.Range("A1").Value = $newText

# --- tasas: update the rate cells ---
$wsTasas = $wb.Worksheets.Item("tasas")

$wsTasas.Range("N10").Value = 425.5
$wsTasas.Range("O10").Value = 3784.62
$wsTasas.Range("N12").Value = 3799.99
$wsTasas.Range("O12").Value = 400.099
